# The sheet holds one weekly price record per row (rows 2..119). This
# edit adds one new weekly record for "Apio" right above the record that
# is currently on row 83 (date 2021-08-08 / serial 44389), pushing that
# row (and everything below it) down by one row. The new record re-uses
# the same min/max/avg price, unit, origin, etc. as the record it now sits
# above, but with its own date and volume:
#   D83 (Fecha)   = 44466
#   J83 (Volumen) = 160
#
# Net effect: dimension grows from A1:R119 to A1:R120, and old rows
# 83..119 become new rows 84..120 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at 83; this shifts rows 83:119 down to 84:120
# (values, formatting/styles and everything else move with the rows).
$ws.Rows("83:83").Insert()

# Seed the freshly-inserted row 83 with the same record that is now on
# row 84 (i.e. what used to be row 83 before the insert), then overwrite
# just the two fields that differ for the new weekly entry.
$src = $ws.Range("A84:R84").Value2
$ws.Range("A83:R83").Value2 = $src

$ws.Range("D83").Value = 44466
$ws.Range("J83").Value = 160
